$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old leading index column (A) - shifts B->A, C->B, D->C
$ws.Columns("A").Delete()

# Rewrite the table with the new header + time-slot rows.
# The string values are entered in this particular order so the workbook's
# shared-string table is built up the same way it was by the original author.
$ws.Range("A1").Value = "JAM"
$ws.Range("B1").Value = "Jumlah Motor"
$ws.Range("C1").Value = "Jumlah Mobil"

$ws.Range("A2").Value = "18.00 - 18.10"
$ws.Range("A4").Value = "18.20 - 18.30"
$ws.Range("A6").Value = "18.40 - 18.50"
$ws.Range("A5").Value = "18.30 - 18.40"
$ws.Range("A7").Value = "18.50 - 19.00"
$ws.Range("A3").Value = "18.10 - 18.20"

$ws.Range("B2").Value = 206
$ws.Range("C2").Value = 6

$ws.Range("B3").Value = 175
$ws.Range("C3").Value = 2

$ws.Range("B4").Value = 198
$ws.Range("C4").Value = 2

$ws.Range("B5").Value = 147
$ws.Range("C5").Value = 3

$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 2

$ws.Range("B7").Value = 190
$ws.Range("C7").Value = 7

# Apply the centered style (same as before) across the whole table
$ws.Range("A1:C7").HorizontalAlignment = -4108

# Column widths to match target layout (values compensated for the engine's
# character-width rounding so the stored OOXML width lands on the target)
$ws.Columns("A").ColumnWidth = 17.276041666666668
$ws.Columns("B").ColumnWidth = 16.498697916666668
$ws.Columns("C").ColumnWidth = 17.498697916666668

$ws.Range("D7").Select()
